# Horarios actualizados Linea 141 - 376
# Refresh the scraped schedule data (new scrape timestamp 08:36:20) across
# the three sheets: LP1912, LP1912-215 and 6203-6173.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:36:20"
$ws.Cells.Item(3,1).Value = "Total filas: 104"
$ws.Cells.Item(71,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(72,3).Value = "215B_EL PATO"
$ws.Cells.Item(77,1).Value = "08:36:20"
$ws.Cells.Item(77,2).Value = "08:38"
$ws.Cells.Item(77,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(77,4).Value = 2
$ws.Cells.Item(78,1).Value = "07:58:19"
$ws.Cells.Item(78,3).Value = "16_SANTA ANA"
$ws.Cells.Item(78,4).Value = 44
$ws.Cells.Item(79,1).Value = "06:57:11"
$ws.Cells.Item(79,2).Value = "08:42"
$ws.Cells.Item(79,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(79,4).Value = 105
$ws.Cells.Item(80,1).Value = "07:19:37"
$ws.Cells.Item(80,3).Value = "14_ABASTO"
$ws.Cells.Item(80,4).Value = 84
$ws.Cells.Item(81,2).Value = "08:43"
$ws.Cells.Item(81,3).Value = "16_SANTA ANA"
$ws.Cells.Item(81,4).Value = 24
$ws.Cells.Item(82,1).Value = "08:36:20"
$ws.Cells.Item(82,2).Value = "08:46"
$ws.Cells.Item(82,3).Value = "16_SANTA ANA"
$ws.Cells.Item(82,4).Value = 10
$ws.Cells.Item(83,1).Value = "08:19:33"
$ws.Cells.Item(83,2).Value = "08:53"
$ws.Cells.Item(83,3).Value = "10_OLMOS"
$ws.Cells.Item(83,4).Value = 34
$ws.Cells.Item(84,1).Value = "06:57:11"
$ws.Cells.Item(84,2).Value = "08:54"
$ws.Cells.Item(84,3).Value = "17_ROMERO"
$ws.Cells.Item(84,4).Value = 117
$ws.Cells.Item(85,1).Value = "07:19:37"
$ws.Cells.Item(85,2).Value = "09:01"
$ws.Cells.Item(85,3).Value = "215A_EL PATO"
$ws.Cells.Item(85,4).Value = 102
$ws.Cells.Item(86,1).Value = "08:19:33"
$ws.Cells.Item(86,2).Value = "09:02"
$ws.Cells.Item(86,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(86,4).Value = 43
$ws.Cells.Item(87,1).Value = "07:45:49"
$ws.Cells.Item(87,2).Value = "09:03"
$ws.Cells.Item(87,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(87,4).Value = 78
$ws.Cells.Item(88,1).Value = "08:36:20"
$ws.Cells.Item(88,2).Value = "09:05"
$ws.Cells.Item(88,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(88,4).Value = 29
$ws.Cells.Item(89,1).Value = "07:19:37"
$ws.Cells.Item(89,2).Value = "09:10"
$ws.Cells.Item(89,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(89,4).Value = 111
$ws.Cells.Item(90,1).Value = "08:36:20"
$ws.Cells.Item(90,2).Value = "09:11"
$ws.Cells.Item(90,3).Value = "16_SANTA ANA"
$ws.Cells.Item(90,4).Value = 35
$ws.Cells.Item(91,1).Value = "08:36:20"
$ws.Cells.Item(91,2).Value = "09:13"
$ws.Cells.Item(91,3).Value = "10_OLMOS"
$ws.Cells.Item(91,4).Value = 37
$ws.Cells.Item(92,1).Value = "07:19:37"
$ws.Cells.Item(92,2).Value = "09:16"
$ws.Cells.Item(92,3).Value = "27_EL RETIRO"
$ws.Cells.Item(92,4).Value = 117
$ws.Cells.Item(93,1).Value = "07:58:19"
$ws.Cells.Item(93,2).Value = "09:17"
$ws.Cells.Item(93,3).Value = "27_EL RETIRO"
$ws.Cells.Item(93,4).Value = 79
$ws.Cells.Item(94,2).Value = "09:21"
$ws.Cells.Item(94,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(94,4).Value = 96
$ws.Cells.Item(95,2).Value = "09:22"
$ws.Cells.Item(95,3).Value = "17_ROMERO"
$ws.Cells.Item(95,4).Value = 97
$ws.Cells.Item(96,2).Value = "09:23"
$ws.Cells.Item(96,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(96,4).Value = 98
$ws.Cells.Item(97,2).Value = "09:23"
$ws.Cells.Item(97,3).Value = "17_ROMERO"
$ws.Cells.Item(97,4).Value = 85
$ws.Cells.Item(98,1).Value = "08:19:33"
$ws.Cells.Item(98,2).Value = "09:25"
$ws.Cells.Item(98,3).Value = "16_SANTA ANA"
$ws.Cells.Item(98,4).Value = 66
$ws.Cells.Item(99,1).Value = "07:45:49"
$ws.Cells.Item(99,2).Value = "09:32"
$ws.Cells.Item(99,4).Value = 107
$ws.Cells.Item(100,1).Value = "07:45:49"
$ws.Cells.Item(100,2).Value = "09:33"
$ws.Cells.Item(100,3).Value = "10_OLMOS"
$ws.Cells.Item(100,4).Value = 108
$ws.Cells.Item(101,1).Value = "07:45:49"
$ws.Cells.Item(101,2).Value = "09:41"
$ws.Cells.Item(101,3).Value = "215C_EL PATO"
$ws.Cells.Item(101,4).Value = 116
$ws.Cells.Item(102,1).Value = "07:58:19"
$ws.Cells.Item(102,2).Value = "09:42"
$ws.Cells.Item(102,3).Value = "215C_EL PATO"
$ws.Cells.Item(102,4).Value = 104
$ws.Cells.Item(102,5).Value = "LP1912"
$ws.Cells.Item(103,1).Value = "07:58:19"
$ws.Cells.Item(103,2).Value = "09:43"
$ws.Cells.Item(103,3).Value = "14_ABASTO"
$ws.Cells.Item(103,4).Value = 105
$ws.Cells.Item(103,5).Value = "LP1912"
$ws.Cells.Item(104,1).Value = "07:58:19"
$ws.Cells.Item(104,2).Value = "09:52"
$ws.Cells.Item(104,3).Value = "15_ABASTO"
$ws.Cells.Item(104,4).Value = 114
$ws.Cells.Item(104,5).Value = "LP1912"
$ws.Cells.Item(105,1).Value = "08:19:33"
$ws.Cells.Item(105,2).Value = "10:10"
$ws.Cells.Item(105,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(105,4).Value = 111
$ws.Cells.Item(105,5).Value = "LP1912"
$ws.Cells.Item(106,1).Value = "08:19:33"
$ws.Cells.Item(106,2).Value = "10:12"
$ws.Cells.Item(106,3).Value = "15_ABASTO"
$ws.Cells.Item(106,4).Value = 113
$ws.Cells.Item(106,5).Value = "LP1912"
$ws.Cells.Item(107,1).Value = "08:36:20"
$ws.Cells.Item(107,2).Value = "10:21"
$ws.Cells.Item(107,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(107,4).Value = 105
$ws.Cells.Item(107,5).Value = "LP1912"
$ws.Cells.Item(108,1).Value = "08:36:20"
$ws.Cells.Item(108,2).Value = "10:22"
$ws.Cells.Item(108,3).Value = "17_ROMERO"
$ws.Cells.Item(108,4).Value = 106
$ws.Cells.Item(108,5).Value = "LP1912"
$ws.Cells.Item(109,1).Value = "08:36:20"
$ws.Cells.Item(109,2).Value = "10:26"
$ws.Cells.Item(109,3).Value = "215A_EL PATO"
$ws.Cells.Item(109,4).Value = 110
$ws.Cells.Item(109,5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:36:20"
$ws.Cells.Item(3,1).Value = "Total filas: 15"
$ws.Cells.Item(20,1).Value = "08:36:20"
$ws.Cells.Item(20,2).Value = "10:26"
$ws.Cells.Item(20,3).Value = "215A_EL PATO"
$ws.Cells.Item(20,4).Value = 110
$ws.Cells.Item(20,5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:36:20"
$ws.Cells.Item(3,1).Value = "Total filas: 22"
$ws.Cells.Item(26,1).Value = "08:36:20"
$ws.Cells.Item(26,2).Value = "09:10"
$ws.Cells.Item(26,3).Value = "215D_LA PLATA"
$ws.Cells.Item(26,4).Value = 34
$ws.Cells.Item(26,5).Value = "L6203"
$ws.Cells.Item(27,1).Value = "08:19:33"
$ws.Cells.Item(27,2).Value = "10:03"
$ws.Cells.Item(27,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(27,4).Value = 104
$ws.Cells.Item(27,5).Value = "L6173"

